# Update column G ("K" - strikeouts) values on the active sheet to the
# regenerated values (K instead of Strike#).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 1
    4  = 2
    6  = 2
    7  = 0
    8  = 1
    9  = 6
    10 = 3
    11 = 2
    12 = 0
    13 = 1
    14 = 2
    15 = 3
    16 = 5
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 2
    26 = 1
    27 = 1
    28 = 0
    29 = 2
    30 = 1
    31 = 2
    32 = 1
    33 = 2
    34 = 1
    35 = 1
    36 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
